$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 19611.875
$ws.Range("I2").Value = 2700
$ws.Range("K2").Value = 2700
$ws.Range("M2").Value = -2587

$ws.Range("H38").Value = 3547.25
$ws.Range("J38").Value = 4663
$ws.Range("L38").Value = 13989
$ws.Range("N38").Value = -14733

$ws.Range("H70").Value = 1246
$ws.Range("J70").Value = 1174.5
$ws.Range("L70").Value = 3523.5
$ws.Range("N70").Value = -4063.5

$ws.Range("H73").Value = 1246
$ws.Range("J73").Value = 1174.5
$ws.Range("L73").Value = 3523.5
$ws.Range("N73").Value = -5395.5

$ws.Range("H80").Value = 776
$ws.Range("I80").Value = 619
$ws.Range("K80").Value = 1857
$ws.Range("M80").Value = -859

$ws.Range("H83").Value = 776
$ws.Range("I83").Value = 619
$ws.Range("K83").Value = 5571
$ws.Range("M83").Value = -579

$ws.Range("H88").Value = 1599
$ws.Range("I88").Value = 1698
$ws.Range("J88").Value = 1500
$ws.Range("K88").Value = 1698
$ws.Range("L88").Value = 1500
$ws.Range("M88").Value = -1292
$ws.Range("N88").Value = -2312

$ws.Range("H91").Value = 1599
$ws.Range("I91").Value = 1698
$ws.Range("J91").Value = 1500
$ws.Range("K91").Value = 1698
$ws.Range("L91").Value = 1500
$ws.Range("M91").Value = -294
$ws.Range("N91").Value = -4308

$ws.Range("H111").Value = 2265.5
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws.Range("H112").Value = 1001.26666
$ws.Range("I112").Value = 750
$ws.Range("J112").Value = 1092.6364
$ws.Range("K112").Value = 2250
$ws.Range("L112").Value = 3277.9092
$ws.Range("M112").Value = -1142
$ws.Range("N112").Value = -5493.9092

$ws.Range("H116").Value = 4528.231
$ws.Range("I116").Value = 4218.6665
$ws.Range("K116").Value = 4218.6665
$ws.Range("M116").Value = -776.6665000000003

$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("N125").ClearContents()

$ws.Range("H132").Value = 1104.7391
$ws.Range("I132").Value = 1113.7727
$ws.Range("J132").Value = 906
$ws.Range("K132").Value = 3341.3181
$ws.Range("L132").Value = 2718
$ws.Range("M132").Value = -811.3181
$ws.Range("N132").Value = -7778

$ws.Range("H135").Value = 4816.1113
$ws.Range("I135").Value = 5787
$ws.Range("J135").Value = 1418
$ws.Range("K135").Value = 52083
$ws.Range("L135").Value = 12762
$ws.Range("M135").Value = -49548
$ws.Range("N135").Value = -17832

$ws.Range("H137").Value = 2179.8125
$ws.Range("I137").Value = 1443.6364
$ws.Range("K137").Value = 4330.9092
$ws.Range("M137").Value = -1780.9092

$ws.Range("H138").Value = 7704.3105
$ws.Range("I138").Value = 4497
$ws.Range("K138").Value = 13491
$ws.Range("M138").Value = -8351

$ws.Range("H141").Value = 5665.923
$ws.Range("I141").Value = 4466.8
$ws.Range("K141").Value = 13400.4
$ws.Range("M141").Value = -8220.400000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2086.6
$ws.Range("I74").Value = 1108.25
$ws.Range("K74").Value = 1108.25
$ws.Range("M74").Value = -234.25

$ws.Range("H77").Value = 2086.6
$ws.Range("I77").Value = 1108.25
$ws.Range("K77").Value = 5541.25
$ws.Range("M77").Value = -1173.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 124966.664
$ws.Range("J132").Value = 124966.664
$ws.Range("L132").Value = 124966.664
$ws.Range("N132").Value = -135086.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2420.5
$ws.Range("J58").Value = 2395
$ws.Range("L58").Value = 2395
$ws.Range("N58").Value = -2801

$ws.Range("H136").Value = 2420.5
$ws.Range("J136").Value = 2395
$ws.Range("L136").Value = 7185
$ws.Range("N136").Value = -12285

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 82370460
$ws.Range("I4").Value = 15000826
$ws.Range("K4").Value = 45002478
$ws.Range("M4").Value = -45002366

$ws.Range("H74").Value = 7499.952
$ws.Range("J74").Value = 7499.952
$ws.Range("L74").Value = 22499.856
$ws.Range("N74").Value = -24621.856

$ws.Range("H77").Value = 7499.952
$ws.Range("J77").Value = 7499.952
$ws.Range("L77").Value = 67499.568
$ws.Range("N77").Value = -78107.568

$ws.Range("H114").Value = 2672.3333
$ws.Range("I114").Value = 2578.4
$ws.Range("J114").Value = 2789.75
$ws.Range("K114").Value = 7735.200000000001
$ws.Range("L114").Value = 8369.25
$ws.Range("M114").Value = -4481.200000000001
$ws.Range("N114").Value = -14877.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5422.727
$ws.Range("I132").Value = 5244.5713
$ws.Range("J132").Value = 5734.5
$ws.Range("K132").Value = 15733.7139
$ws.Range("L132").Value = 17203.5
$ws.Range("M132").Value = -13203.7139
$ws.Range("N132").Value = -22263.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws.Range("H132").Value = 5012.5713
$ws.Range("I132").Value = 2549.5
$ws.Range("K132").Value = 7648.5
$ws.Range("M132").Value = -5118.5

$ws.Range("H136").Value = 24462.137
$ws.Range("I136").Value = 3022.125
$ws.Range("K136").Value = 9066.375
$ws.Range("M136").Value = -6516.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 366.5
$ws.Range("I107").Value = 340.2
$ws.Range("K107").Value = 1020.6
$ws.Range("M107").Value = 899.4000000000001

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws.Range("H132").Value = 2506.1516
$ws.Range("J132").Value = 3348
$ws.Range("L132").Value = 10044
$ws.Range("N132").Value = -15104

$ws.Range("H136").Value = 10936.357
$ws.Range("I136").Value = 12313.909
$ws.Range("J136").Value = 5885.3335
$ws.Range("K136").Value = 36941.727
$ws.Range("L136").Value = 17656.0005
$ws.Range("M136").Value = -34391.727
$ws.Range("N136").Value = -22756.0005
